# TokenIteratorFieldRewriterSplit update:
# The template used to embed the M2Doc expression
#     m:'doc.html'.fromHTMLURI()
# as a real Word field (fldChar begin / instrText* / fldChar end) wrapped
# around the `_GoBack` bookmark. The new parser instead expects the plain
# M2Doc token delimiters `{ ... }` spelled out as literal run text (still
# split across several runs, still straddling the same bookmark), so the
# field machinery is replaced by plain <w:t> runs reading:
#   { m : ' doc.html [bookmark] '.fromHTMLURI() }
$d = $word.ActiveDocument

# Locate the paragraph that carries the M2Doc field (there is exactly one
# in this document).
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Fields.Count -gt 0) {
        $targetPara = $p
    }
}

if ($targetPara -ne $null) {
    $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00C52979" w:rsidRDefault="00C52979" w:rsidP="00F5495F">' `
      + '<w:r><w:t>{</w:t></w:r>' `
      + '<w:r><w:t>m</w:t></w:r>' `
      + '<w:r><w:t>:</w:t></w:r>' `
      + "<w:r><w:t>'</w:t></w:r>" `
      + '<w:r><w:t>doc.html</w:t></w:r>' `
      + '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' `
      + "<w:r><w:t>'.fromHTMLURI()</w:t></w:r>" `
      + '<w:r><w:t xml:space="preserve">}</w:t></w:r>' `
      + '</w:p>'

    $targetPara.Range.InsertXML($xml) | Out-Null
}
